$d = $word.ActiveDocument

# 1. Remove the paragraph that holds the old "Meta description" label and
#    its text (it directly follows the document title at the very top).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Meta description*") {
        $p.Range.Delete()
        break
    }
}

# 2. Find the paragraph that still holds the old image-prompt text and
#    insert a new paragraph with the (now bold, non-heading) title text
#    right before it.
$imageOld = "Create a feature image fitting the game Cubes"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$imageOld*") {
        $prevPara = $d.Paragraphs.Item($i - 1)
        $prevPara.Range.InsertParagraphAfter()

        $titlePara = $d.Paragraphs.Item($i)
        $titlePara.Style = "Normal"
        $titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
        $titleRange.Text = "Play Cubes Slot Game for Free | Cluster Gameplay and Unique Wins"
        $titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End - 1)
        $titleRange.Bold = 1
        break
    }
}

# 3. Replace the old image-prompt paragraph's text with the meta
#    description text, keeping the paragraph's italic run formatting.
$old = "Create a feature image fitting the game Cubes: - Style: Cartoon - Subject: A happy Maya warrior with glasses as the main focus of the image - Use bright and bold colors to capture the fun and playful nature of the game - Include elements of the Rubik's Cube to tie in the game's theme and give a nod to its inspiration - The warrior can be shown holding a Rubik's Cube or standing next to a larger-than-life version of the cube - The image should showcase the lively and entertaining nature of the game and inspire players to try their luck spinning the colorful cubes."
$new = "Find out about Cubes, a cluster game with an expanding grid and exciting special features such as Central Multiplier and Free Spins. Play the game for free now."
[void]$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
